$wb = $excel.ActiveWorkbook

# Rename the "Id" column header to "Id Proceso" on the Procesos sheet
$ws = $wb.Worksheets.Item("Procesos")
$ws.Range("A1").Value = "Id Proceso"

# Widen column A to fit the new, longer header text (no longer auto/best-fit)
$ws.Columns.Item(1).ColumnWidth = 19
